$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 79 by copying formatting from row 78 (A78:D78 -> A79:D79);
# this also extends the sheet dimension to D79 automatically.
$ws.Range("A78:D78").Copy()
$ws.Range("A79:D79").PasteSpecial(-4122)

# Re-populate rows 2-79 with the newly sorted error-relative data
$ws.Range("A2").Value = 40
$ws.Range("B2").Value = 2233000
$ws.Range("C2").Value = 4321220
$ws.Range("D2").Value = 93.51634572324228
$ws.Range("A3").Value = 25
$ws.Range("B3").Value = 1750000
$ws.Range("C3").Value = 3380407.25
$ws.Range("D3").Value = 93.16612857142857
$ws.Range("A4").Value = 61
$ws.Range("B4").Value = 1750000
$ws.Range("C4").Value = 3272603
$ws.Range("D4").Value = 87.00588571428571
$ws.Range("A5").Value = 55
$ws.Range("B5").Value = 2380000
$ws.Range("C5").Value = 4371614.5
$ws.Range("D5").Value = 83.68128151260504
$ws.Range("A6").Value = 9
$ws.Range("B6").Value = 3500000
$ws.Range("C6").Value = 5932616
$ws.Range("D6").Value = 69.5033142857143
$ws.Range("A7").Value = 65
$ws.Range("B7").Value = 1750000
$ws.Range("C7").Value = 2822316
$ws.Range("D7").Value = 61.2752
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 2275000
$ws.Range("C8").Value = 3535271.25
$ws.Range("D8").Value = 55.39653846153846
$ws.Range("A9").Value = 66
$ws.Range("B9").Value = 8400000
$ws.Range("C9").Value = 4119406.75
$ws.Range("D9").Value = 50.95944345238095
$ws.Range("A10").Value = 77
$ws.Range("B10").Value = 2520000
$ws.Range("C10").Value = 3699640.75
$ws.Range("D10").Value = 46.81114087301587
$ws.Range("A11").Value = 14
$ws.Range("B11").Value = 4025000
$ws.Range("C11").Value = 5738945.5
$ws.Range("D11").Value = 42.58249689440994
$ws.Range("A12").Value = 71
$ws.Range("B12").Value = 4340000
$ws.Range("C12").Value = 6158072.5
$ws.Range("D12").Value = 41.89107142857143
$ws.Range("A13").Value = 60
$ws.Range("B13").Value = 3570000
$ws.Range("C13").Value = 5003962.5
$ws.Range("D13").Value = 40.16701680672267
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = 1890000
$ws.Range("C14").Value = 2630729.75
$ws.Range("D14").Value = 39.19205026455028
$ws.Range("A15").Value = 58
$ws.Range("B15").Value = 2450000
$ws.Range("C15").Value = 3392808.75
$ws.Range("D15").Value = 38.48198979591837
$ws.Range("A16").Value = 23
$ws.Range("B16").Value = 4200000
$ws.Range("C16").Value = 5768288
$ws.Range("D16").Value = 37.34019047619048
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = 3010000
$ws.Range("C17").Value = 4128174.5
$ws.Range("D17").Value = 37.14865448504983
$ws.Range("A18").Value = 50
$ws.Range("B18").Value = 2450000
$ws.Range("C18").Value = 3340312.5
$ws.Range("D18").Value = 36.33928571428572
$ws.Range("A19").Value = 31
$ws.Range("B19").Value = 5110000
$ws.Range("C19").Value = 3332437.25
$ws.Range("D19").Value = 34.78596379647749
$ws.Range("A20").Value = 45
$ws.Range("B20").Value = 3080000
$ws.Range("C20").Value = 4132898
$ws.Range("D20").Value = 34.185
$ws.Range("A21").Value = 38
$ws.Range("B21").Value = 3850000
$ws.Range("C21").Value = 5140917.5
$ws.Range("D21").Value = 33.53032467532467
$ws.Range("A22").Value = 29
$ws.Range("B22").Value = 2135000
$ws.Range("C22").Value = 2782056
$ws.Range("D22").Value = 30.30707259953162
$ws.Range("A23").Value = 30
$ws.Range("B23").Value = 3500000
$ws.Range("C23").Value = 4558037
$ws.Range("D23").Value = 30.22962857142857
$ws.Range("A24").Value = 57
$ws.Range("B24").Value = 2870000
$ws.Range("C24").Value = 3716983.5
$ws.Range("D24").Value = 29.51162020905923
$ws.Range("A25").Value = 4
$ws.Range("B25").Value = 3290000
$ws.Range("C25").Value = 4245402
$ws.Range("D25").Value = 29.03957446808511
$ws.Range("A26").Value = 42
$ws.Range("B26").Value = 8680000
$ws.Range("C26").Value = 6177615
$ws.Range("D26").Value = 28.82932027649769
$ws.Range("A27").Value = 33
$ws.Range("B27").Value = 3430000
$ws.Range("C27").Value = 4388945.5
$ws.Range("D27").Value = 27.95759475218659
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = 4543000
$ws.Range("C28").Value = 5811781.5
$ws.Range("D28").Value = 27.92827426810478
$ws.Range("A29").Value = 10
$ws.Range("B29").Value = 3080000
$ws.Range("C29").Value = 3918247.75
$ws.Range("D29").Value = 27.21583603896104
$ws.Range("A30").Value = 41
$ws.Range("B30").Value = 5600000
$ws.Range("C30").Value = 4100026.25
$ws.Range("D30").Value = 26.78524553571429
$ws.Range("A31").Value = 53
$ws.Range("B31").Value = 2870000
$ws.Range("C31").Value = 3598649.5
$ws.Range("D31").Value = 25.38848432055749
$ws.Range("A32").Value = 26
$ws.Range("B32").Value = 4375000
$ws.Range("C32").Value = 3283094.5
$ws.Range("D32").Value = 24.95784
$ws.Range("A33").Value = 36
$ws.Range("B33").Value = 3990000
$ws.Range("C33").Value = 4984787
$ws.Range("D33").Value = 24.93200501253133
$ws.Range("A34").Value = 24
$ws.Range("B34").Value = 3640000
$ws.Range("C34").Value = 4508519.5
$ws.Range("D34").Value = 23.86042582417581
$ws.Range("A35").Value = 63
$ws.Range("B35").Value = 5033000
$ws.Range("C35").Value = 3866840.75
$ws.Range("D35").Value = 23.17026127558116
$ws.Range("A36").Value = 19
$ws.Range("B36").Value = 2485000
$ws.Range("C36").Value = 3055258.25
$ws.Range("D36").Value = 22.94801810865191
$ws.Range("A37").Value = 44
$ws.Range("B37").Value = 8294999.999999999
$ws.Range("C37").Value = 6462561.5
$ws.Range("D37").Value = 22.09088004822181
$ws.Range("A38").Value = 69
$ws.Range("B38").Value = 3703000
$ws.Range("C38").Value = 4501734.5
$ws.Range("D38").Value = 21.56992978665948
$ws.Range("A39").Value = 54
$ws.Range("B39").Value = 3290000
$ws.Range("C39").Value = 3967805.25
$ws.Range("D39").Value = 20.60198328267477
$ws.Range("A40").Value = 43
$ws.Range("B40").Value = 3640000
$ws.Range("C40").Value = 4362461.5
$ws.Range("D40").Value = 19.84784340659339
$ws.Range("A41").Value = 17
$ws.Range("B41").Value = 3710000
$ws.Range("C41").Value = 2978359.75
$ws.Range("D41").Value = 19.7207614555256
$ws.Range("A42").Value = 32
$ws.Range("B42").Value = 5530000
$ws.Range("C42").Value = 4461440.5
$ws.Range("D42").Value = 19.32295660036166
$ws.Range("A43").Value = 72
$ws.Range("B43").Value = 5950000
$ws.Range("C43").Value = 4830545
$ws.Range("D43").Value = 18.81436974789916
$ws.Range("A44").Value = 51
$ws.Range("B44").Value = 5565000
$ws.Range("C44").Value = 6569257.5
$ws.Range("D44").Value = 18.04595687331537
$ws.Range("A45").Value = 0
$ws.Range("B45").Value = 4690000
$ws.Range("C45").Value = 3860289.5
$ws.Range("D45").Value = 17.69105543710021
$ws.Range("A46").Value = 68
$ws.Range("B46").Value = 5215000
$ws.Range("C46").Value = 4305643.5
$ws.Range("D46").Value = 17.43732502396932
$ws.Range("A47").Value = 6
$ws.Range("B47").Value = 5495000
$ws.Range("C47").Value = 4548168.5
$ws.Range("D47").Value = 17.23078252957234
$ws.Range("A48").Value = 13
$ws.Range("B48").Value = 7245000
$ws.Range("C48").Value = 8487414
$ws.Range("D48").Value = 17.14857142857143
$ws.Range("A49").Value = 34
$ws.Range("B49").Value = 2408000
$ws.Range("C49").Value = 2816920.25
$ws.Range("D49").Value = 16.98173795681063
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 6300000
$ws.Range("C50").Value = 5252373
$ws.Range("D50").Value = 16.629
$ws.Range("A51").Value = 75
$ws.Range("B51").Value = 6650000
$ws.Range("C51").Value = 5554134
$ws.Range("D51").Value = 16.47918796992481
$ws.Range("A52").Value = 2
$ws.Range("B52").Value = 4060000
$ws.Range("C52").Value = 3407564.5
$ws.Range("D52").Value = 16.06983990147783
$ws.Range("A53").Value = 28
$ws.Range("B53").Value = 5390000
$ws.Range("C53").Value = 6251707
$ws.Range("D53").Value = 15.98714285714286
$ws.Range("A54").Value = 21
$ws.Range("B54").Value = 6419000
$ws.Range("C54").Value = 5408772.5
$ws.Range("D54").Value = 15.73808225580308
$ws.Range("A55").Value = 15
$ws.Range("B55").Value = 3430000
$ws.Range("C55").Value = 2903050.5
$ws.Range("D55").Value = 15.36295918367347
$ws.Range("A56").Value = 47
$ws.Range("B56").Value = 3885000
$ws.Range("C56").Value = 4466045
$ws.Range("D56").Value = 14.95611325611326
$ws.Range("A57").Value = 18
$ws.Range("B57").Value = 6930000
$ws.Range("C57").Value = 5912520.5
$ws.Range("D57").Value = 14.68224386724387
$ws.Range("A58").Value = 62
$ws.Range("B58").Value = 3465000
$ws.Range("C58").Value = 3964181.25
$ws.Range("D58").Value = 14.40638528138528
$ws.Range("A59").Value = 20
$ws.Range("B59").Value = 2835000
$ws.Range("C59").Value = 3240659
$ws.Range("D59").Value = 14.3089594356261
$ws.Range("A60").Value = 64
$ws.Range("B60").Value = 4907000
$ws.Range("C60").Value = 4205315.5
$ws.Range("D60").Value = 14.29966374566945
$ws.Range("A61").Value = 35
$ws.Range("B61").Value = 3087000
$ws.Range("C61").Value = 3522897.75
$ws.Range("D61").Value = 14.12043245869777
$ws.Range("A62").Value = 48
$ws.Range("B62").Value = 5250000
$ws.Range("C62").Value = 5983282
$ws.Range("D62").Value = 13.96727619047619
$ws.Range("A63").Value = 12
$ws.Range("B63").Value = 3920000
$ws.Range("C63").Value = 4450266
$ws.Range("D63").Value = 13.52719387755102
$ws.Range("A64").Value = 37
$ws.Range("B64").Value = 5810000
$ws.Range("C64").Value = 5033414
$ws.Range("D64").Value = 13.36636833046471
$ws.Range("A65").Value = 7
$ws.Range("B65").Value = 4200000
$ws.Range("C65").Value = 3638951.25
$ws.Range("D65").Value = 13.35830357142857
$ws.Range("A66").Value = 27
$ws.Range("B66").Value = 2590000
$ws.Range("C66").Value = 2933935.25
$ws.Range("D66").Value = 13.27935328185328
$ws.Range("A67").Value = 73
$ws.Range("B67").Value = 4060000
$ws.Range("C67").Value = 4577985.5
$ws.Range("D67").Value = 12.75826354679803
$ws.Range("A68").Value = 70
$ws.Range("B68").Value = 2940000
$ws.Range("C68").Value = 3297328.75
$ws.Range("D68").Value = 12.15403911564626
$ws.Range("A69").Value = 8
$ws.Range("B69").Value = 5460000
$ws.Range("C69").Value = 6122246.5
$ws.Range("D69").Value = 12.12905677655678
$ws.Range("A70").Value = 59
$ws.Range("B70").Value = 5950000
$ws.Range("C70").Value = 5228970
$ws.Range("D70").Value = 12.1181512605042
$ws.Range("A71").Value = 56
$ws.Range("B71").Value = 8400000
$ws.Range("C71").Value = 7413476
$ws.Range("D71").Value = 11.74433333333333
$ws.Range("A72").Value = 16
$ws.Range("B72").Value = 4550000
$ws.Range("C72").Value = 5068953
$ws.Range("D72").Value = 11.40556043956044
$ws.Range("A73").Value = 74
$ws.Range("B73").Value = 3360000
$ws.Range("C73").Value = 2979431
$ws.Range("D73").Value = 11.32645833333333
$ws.Range("A74").Value = 39
$ws.Range("B74").Value = 3780000
$ws.Range("C74").Value = 4203505.5
$ws.Range("D74").Value = 11.20384920634922
$ws.Range("A75").Value = 67
$ws.Range("B75").Value = 3675000
$ws.Range("C75").Value = 4085956.25
$ws.Range("D75").Value = 11.18248299319729
$ws.Range("A76").Value = 52
$ws.Range("B76").Value = 4970000
$ws.Range("C76").Value = 5518557.5
$ws.Range("D76").Value = 11.03737424547284
$ws.Range("A77").Value = 22
$ws.Range("B77").Value = 3675000
$ws.Range("C77").Value = 3271173
$ws.Range("D77").Value = 10.98848979591836
$ws.Range("A78").Value = 76
$ws.Range("B78").Value = 7419999.999999999
$ws.Range("C78").Value = 6640457.5
$ws.Range("D78").Value = 10.50596361185983
$ws.Range("A79").Value = 46
$ws.Range("B79").Value = 6090000
$ws.Range("C79").Value = 5455917.5
$ws.Range("D79").Value = 10.41186371100164
